$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1647.6177
$ws.Range("J17").Value = 1764.6333
$ws.Range("L17").Value = 5293.8999
$ws.Range("N17").Value = -5629.8999
$ws.Range("H113").Value = 94784750
$ws.Range("I113").Value = 18520822
$ws.Range("J113").Value = 136383260
$ws.Range("K113").Value = 18520822
$ws.Range("L113").Value = 136383260
$ws.Range("M113").Value = -18517568
$ws.Range("N113").Value = -136389768
$ws.Range("H138").Value = 1432650.6
$ws.Range("J138").Value = 2569743.2
$ws.Range("L138").Value = 7709229.600000001
$ws.Range("N138").Value = -7719509.600000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 26317024
$ws.Range("I2").Value = 986.8929000000001
$ws.Range("K2").Value = 986.8929000000001
$ws.Range("M2").Value = -873.8929000000001
$ws.Range("H34").Value = 161399.6
$ws.Range("I34").Value = 212666.67
$ws.Range("J34").Value = 84499
$ws.Range("K34").Value = 212666.67
$ws.Range("L34").Value = 84499
$ws.Range("M34").Value = -212395.67
$ws.Range("N34").Value = -85041
$ws.Range("H45").Value = 5261.5625
$ws.Range("I45").Value = 1736.1428
$ws.Range("J45").Value = 8003.5557
$ws.Range("K45").Value = 1736.1428
$ws.Range("L45").Value = 8003.5557
$ws.Range("M45").Value = -1359.1428
$ws.Range("N45").Value = -8757.555700000001
$ws.Range("H104").Value = 44944
$ws.Range("J104").Value = 44944
$ws.Range("L104").Value = 44944
$ws.Range("N104").Value = -51932
$ws.Range("H110").Value = 66668084
$ws.Range("I110").Value = 1745
$ws.Range("J110").Value = 111112310
$ws.Range("K110").Value = 1745
$ws.Range("L110").Value = 111112310
$ws.Range("M110").Value = 300
$ws.Range("N110").Value = -111116400
$ws.Range("H116").Value = 26317024
$ws.Range("I116").Value = 986.8929000000001
$ws.Range("K116").Value = 986.8929000000001
$ws.Range("M116").Value = 1307.1071
$ws.Range("H132").Value = 6322.606
$ws.Range("I132").Value = 2225.1333
$ws.Range("K132").Value = 6675.3999
$ws.Range("M132").Value = -4145.3999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 26317024
$ws.Range("I3").Value = 986.8929000000001
$ws.Range("K3").Value = 986.8929000000001
$ws.Range("M3").Value = -872.8929000000001
$ws.Range("H22").Value = 7936893.5
$ws.Range("I22").Value = 9259626
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 9259626
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -9259453
$ws.Range("N22").Value = -846
$ws.Range("H94").Value = 3146.4736
$ws.Range("I94").Value = 1438.5333
$ws.Range("K94").Value = 1438.5333
$ws.Range("M94").Value = -987.5333000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7977.3076
$ws.Range("I31").Value = 2776.625
$ws.Range("J31").Value = 11595.174
$ws.Range("K31").Value = 2776.625
$ws.Range("L31").Value = 11595.174
$ws.Range("M31").Value = -2481.625
$ws.Range("N31").Value = -12185.174
$ws.Range("H34").Value = 7977.3076
$ws.Range("I34").Value = 2776.625
$ws.Range("J34").Value = 11595.174
$ws.Range("K34").Value = 2776.625
$ws.Range("L34").Value = 11595.174
$ws.Range("M34").Value = -2574.625
$ws.Range("N34").Value = -11999.174

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1770.5
$ws.Range("I5").Value = 960.36365
$ws.Range("J5").Value = 3043.5715
$ws.Range("K5").Value = 2881.09095
$ws.Range("L5").Value = 9130.7145
$ws.Range("M5").Value = -2769.09095
$ws.Range("N5").Value = -9354.7145
$ws.Range("H14").Value = 20844650
$ws.Range("I14").Value = 20844650
$ws.Range("K14").Value = 62533950
$ws.Range("M14").Value = -62533777
$ws.Range("H105").Value = 10000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 30000
$ws.Range("M105").Value = $null
$ws.Range("N105").Value = -35242
$ws.Range("H108").Value = 1163.4286
$ws.Range("I108").Value = 1163.4286
$ws.Range("K108").Value = 3490.2858
$ws.Range("M108").Value = -610.2857999999997
$ws.Range("H113").Value = 5700.9443
$ws.Range("J113").Value = 9091.299999999999
$ws.Range("L113").Value = 27273.9
$ws.Range("N113").Value = -31613.9
$ws.Range("H119").Value = 3312.6667
$ws.Range("I119").Value = 2747
$ws.Range("K119").Value = 8241
$ws.Range("M119").Value = -3403
$ws.Range("H131").Value = 1642.7142
$ws.Range("I131").Value = 1428.375
$ws.Range("J131").Value = 1728.45
$ws.Range("K131").Value = 4285.125
$ws.Range("L131").Value = 5185.35
$ws.Range("M131").Value = 754.875
$ws.Range("N131").Value = -15265.35
$ws.Range("H135").Value = 1770.5
$ws.Range("I135").Value = 960.36365
$ws.Range("J135").Value = 3043.5715
$ws.Range("K135").Value = 8643.272849999999
$ws.Range("L135").Value = 27392.1435
$ws.Range("M135").Value = -6108.272849999999
$ws.Range("N135").Value = -32462.1435
$ws.Range("H141").Value = 4407.6665
$ws.Range("I141").Value = 4407.6665
$ws.Range("K141").Value = 13222.9995
$ws.Range("M141").Value = -8042.999500000002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 29285.572
$ws.Range("J46").Value = 43749.75
$ws.Range("L46").Value = 43749.75
$ws.Range("N46").Value = -44061.75
$ws.Range("H70").Value = 6012.7144
$ws.Range("I70").Value = 4624.364
$ws.Range("J70").Value = 7539.9
$ws.Range("K70").Value = 4624.364
$ws.Range("L70").Value = 7539.9
$ws.Range("M70").Value = -4354.364
$ws.Range("N70").Value = -8079.9
$ws.Range("H73").Value = 6012.7144
$ws.Range("I73").Value = 4624.364
$ws.Range("J73").Value = 7539.9
$ws.Range("K73").Value = 4624.364
$ws.Range("L73").Value = 7539.9
$ws.Range("M73").Value = -3688.364
$ws.Range("N73").Value = -9411.9
$ws.Range("H107").Value = 882.2593000000001
$ws.Range("I107").Value = 1420.75
$ws.Range("J107").Value = 655.5263
$ws.Range("K107").Value = 1420.75
$ws.Range("L107").Value = 655.5263
$ws.Range("M107").Value = 499.25
$ws.Range("N107").Value = -4495.5263
$ws.Range("H126").Value = 6449.8
$ws.Range("I126").Value = 3067.5557
$ws.Range("K126").Value = 9202.667099999999
$ws.Range("M126").Value = -6732.667099999999
$ws.Range("H132").Value = 5304.654
$ws.Range("I132").Value = 2055.1765
$ws.Range("K132").Value = 6165.529500000001
$ws.Range("M132").Value = -3635.529500000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1504.6364
$ws.Range("I22").Value = 387
$ws.Range("J22").Value = 2143.2856
$ws.Range("K22").Value = 387
$ws.Range("L22").Value = 2143.2856
$ws.Range("M22").Value = -92
$ws.Range("N22").Value = -2733.2856
$ws.Range("H27").Value = 1504.6364
$ws.Range("I27").Value = 387
$ws.Range("J27").Value = 2143.2856
$ws.Range("K27").Value = 387
$ws.Range("L27").Value = 2143.2856
$ws.Range("M27").Value = -280
$ws.Range("N27").Value = -2357.2856
$ws.Range("H46").Value = 1233597.1
$ws.Range("I46").Value = 2156548.8
$ws.Range("K46").Value = 2156548.8
$ws.Range("M46").Value = -2156360.8
$ws.Range("H59").Value = 56134
$ws.Range("J59").Value = 56134
$ws.Range("L59").Value = 56134
$ws.Range("N59").Value = -57442

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 29573.857
$ws.Range("I45").Value = 24664
$ws.Range("J45").Value = 30912.908
$ws.Range("K45").Value = 24664
$ws.Range("L45").Value = 30912.908
$ws.Range("M45").Value = -24173
$ws.Range("N45").Value = -31894.908
$ws.Range("H113").Value = 12659.392
$ws.Range("I113").Value = 24241.273
$ws.Range("K113").Value = 72723.819
$ws.Range("M113").Value = -70553.819
$ws.Range("H130").Value = 48197.145
$ws.Range("J130").Value = 54563.332
$ws.Range("L130").Value = 54563.332
$ws.Range("N130").Value = -64603.332
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null
